$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells must be forced to Text format first so that
# Excel does not auto-convert numeric-looking strings (e.g. "1.00", "0.100")
# into actual numbers and lose their exact textual representation.
$priceCells = @("D2", "D3", "D5", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D25", "D27", "D28", "D30", "D33", "D34", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D50", "D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.185.11'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '2.056.26'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '248.57'
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D8').Value = '57.41'
$ws.Range('E8').Value = '  -1.86%  '
$ws.Range('D9').Value = '0.386'
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('D10').Value = '0.0785'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D12').Value = '16.15'
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').Value = '0.917'
$ws.Range('E13').Value = '  +13.75%  '
$ws.Range('D14').Value = '2.357.21'
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').Value = '5.73'
$ws.Range('D16').Value = '2.057.53'
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('E17').Value = '  +12.47%  '
$ws.Range('D18').Value = '37.276.92'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('D19').Value = '75.06'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = '0.0₃0898'
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('D21').Value = '5.48'
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('D22').Value = '238.05'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  +4.72%  '
$ws.Range('D25').Value = '9.63'
$ws.Range('E25').Value = '  +3.81%  '
$ws.Range('E26').Value = '  -3.52%  '
$ws.Range('D27').Value = '170.72'
$ws.Range('E27').Value = '  +1.04%  '
$ws.Range('D28').Value = '20.24'
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').Value = '5.14'
$ws.Range('E30').Value = '  +8.68%  '
$ws.Range('E31').Value = '  +3.28%  '
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('D33').Value = '4.65'
$ws.Range('E33').Value = '  +3.94%  '
$ws.Range('D34').Value = '0.0882'
$ws.Range('E34').Value = '  +2.10%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  +3.85%  '
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('D39').Value = '5.17'
$ws.Range('E39').Value = '  +11.19%  '
$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').Value = '3.07'
$ws.Range('E40').Value = '  +8.28%  '
$ws.Range('B41').Value = 'Cronos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D41').Value = '0.100'
$ws.Range('E41').Value = '  -9.35%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Value = '17.65'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0224'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').Value = '1.16'
$ws.Range('E44').Value = '  +2.64%  '
$ws.Range('D45').Value = '96.71'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').Value = '2.42'
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('D47').Value = '1.277.69'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('E48').Value = '  -0.90%  '
$ws.Range('E49').Value = '  +1.41%  '
$ws.Range('D50').Value = '2.245.37'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('D51').Value = '0.147'
$ws.Range('E51').Value = '  +9.40%  '

# Restore default (Normal) style on the price cells so no extraneous
# formatting/style differences are introduced versus the original file.
foreach ($cell in $priceCells) {
    $ws.Range($cell).Style = "Normal"
}
